$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$rowCount = $usedRange.Rows.Count

for ($r = 1; $r -le $rowCount; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    $val = $cell.Value()
    if ($val -ne $null -and $val -is [string] -and $val.StartsWith("System, ")) {
        $parts = $val -split ", "
        if ($parts.Count -ge 2) {
            $tmp = $parts[0]
            $parts[0] = $parts[1]
            $parts[1] = $tmp
            $cell.Value = [string]::Join(", ", $parts)
        }
    }
}
